$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.305.81"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.609.64"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").Value = "'212.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").Value = "'0.249"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("D10").Value = "'18.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "1.827.97"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").Value = "1.608.91"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("E14").Value = "  +0.85%  "

$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").Value = "26.288.90"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'61.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "'203.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("D21").Value = "'4.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("D22").Value = "'9.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.02%  "

$ws.Range("D23").Value = "'6.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "'1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.92%  "

$ws.Range("D25").Value = "'144.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("E27").Value = "  -5.64%  "

$ws.Range("D28").Value = "'15.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("D29").Value = "'6.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.67%  "

$ws.Range("D30").Value = "'0.0489"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.85%  "

$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").Value = "'2.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.40%  "

$ws.Range("E34").Value = "  +2.84%  "

$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "1.155.57"
$ws.Range("E36").Value = "  +4.45%  "

$ws.Range("E37").Value = "  +8.71%  "

$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "'0.501"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").Value = "'0.785"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("D44").Value = "1.741.40"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "'91.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("D47").Value = "'54.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").Value = "0.0₇0962"
$ws.Range("E49").Value = "  -11.99%  "

$ws.Range("D50").Value = "'0.406"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
